$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 670.2143
$ws.Range("I2").Value = 721.5
$ws.Range("J2").Value = 631.75
$ws.Range("K2").Value = 721.5
$ws.Range("L2").Value = 631.75
$ws.Range("M2").Value = -608.5
$ws.Range("N2").Value = -857.75
$ws.Range("H9").Value = 5518.05
$ws.Range("I9").Value = 10163.3
$ws.Range("K9").Value = 10163.3
$ws.Range("M9").Value = -9994.299999999999
$ws.Range("H19").Value = 2287
$ws.Range("I19").Value = 1243.5333
$ws.Range("J19").Value = 3265.25
$ws.Range("K19").Value = 1243.5333
$ws.Range("L19").Value = 3265.25
$ws.Range("M19").Value = -1068.5333
$ws.Range("N19").Value = -3615.25
$ws.Range("H43").Value = 4128.4287
$ws.Range("I43").Value = 3749.75
$ws.Range("J43").Value = 4633.3335
$ws.Range("K43").Value = 3749.75
$ws.Range("L43").Value = 4633.3335
$ws.Range("M43").Value = -3680.75
$ws.Range("N43").Value = -4771.3335
$ws.Range("H116").Value = 3850.75
$ws.Range("I116").Value = 5783.3335
$ws.Range("K116").Value = 5783.3335
$ws.Range("M116").Value = -2341.3335
$ws.Range("H125").Value = 6500
$ws.Range("J125").Value = 6500
$ws.Range("L125").Value = 58500
$ws.Range("N125").Value = -63420
$ws.Range("H138").Value = 9808673
$ws.Range("I138").Value = 1422.9286
$ws.Range("K138").Value = 4268.7858
$ws.Range("M138").Value = 871.2142000000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4170.7646
$ws.Range("I32").Value = 2864.7112
$ws.Range("K32").Value = 2864.7112
$ws.Range("M32").Value = -2577.7112
$ws.Range("H45").Value = 2015.2354
$ws.Range("I45").Value = 1875.6154
$ws.Range("K45").Value = 1875.6154
$ws.Range("M45").Value = -1498.6154
$ws.Range("H76").Value = 65999.664
$ws.Range("J76").Value = 65999.664
$ws.Range("L76").Value = 65999.664
$ws.Range("N76").Value = -66675.664
$ws.Range("H79").Value = 65999.664
$ws.Range("J79").Value = 65999.664
$ws.Range("L79").Value = 65999.664
$ws.Range("N79").Value = -68339.664
$ws.Range("H80").Value = 80000
$ws.Range("J80").Value = 80000
$ws.Range("L80").Value = 80000
$ws.Range("N80").Value = -81996
$ws.Range("H83").Value = 80000
$ws.Range("J83").Value = 80000
$ws.Range("L83").Value = 240000
$ws.Range("N83").Value = -249984
$ws.Range("H124").Value = 44699.75
$ws.Range("J124").Value = 44699.75
$ws.Range("L124").Value = 44699.75
$ws.Range("N124").Value = -54519.75
$ws.Range("H132").Value = 40059116
$ws.Range("I132").Value = 3235.1428
$ws.Range("K132").Value = 9705.428400000001
$ws.Range("M132").Value = -7175.428400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 12515.077
$ws.Range("I86").Value = 5330.7393
$ws.Range("K86").Value = 5330.7393
$ws.Range("M86").Value = -4207.7393
$ws.Range("H89").Value = 12515.077
$ws.Range("I89").Value = 5330.7393
$ws.Range("K89").Value = 26653.6965
$ws.Range("M89").Value = -21037.6965
$ws.Range("H110").Value = 51500
$ws.Range("J110").Value = 51500
$ws.Range("L110").Value = 51500
$ws.Range("N110").Value = -59680
$ws.Range("H134").Value = 2976.0386
$ws.Range("I134").Value = 2733.8262
$ws.Range("J134").Value = 4833
$ws.Range("K134").Value = 8201.4786
$ws.Range("L134").Value = 14499
$ws.Range("M134").Value = -5666.4786
$ws.Range("N134").Value = -19569
$ws.Range("H140").Value = 119080
$ws.Range("J140").Value = 119080
$ws.Range("L140").Value = 119080
$ws.Range("N140").Value = -129440

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 34725
$ws.Range("J88").Value = 34725
$ws.Range("L88").Value = 34725
$ws.Range("N88").Value = -35537
$ws.Range("H91").Value = 34725
$ws.Range("J91").Value = 34725
$ws.Range("L91").Value = 34725
$ws.Range("N91").Value = -37533
$ws.Range("H110").Value = 42665
$ws.Range("J110").Value = 53997.5
$ws.Range("L110").Value = 53997.5
$ws.Range("N110").Value = -62177.5
$ws.Range("H116").Value = 55995
$ws.Range("J116").Value = 55995
$ws.Range("L116").Value = 55995
$ws.Range("N116").Value = -65173
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").ClearContents()
$ws.Range("N125").Value = 0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 264.8
$ws.Range("I2").Value = 377.3846
$ws.Range("J2").Value = 55.714287
$ws.Range("K2").Value = 2264.3076
$ws.Range("L2").Value = 334.285722
$ws.Range("M2").Value = -2151.3076
$ws.Range("N2").Value = -560.285722
$ws.Range("H131").Value = 24935.709
$ws.Range("I131").Value = 61984.176
$ws.Range("J131").Value = 4618.8066
$ws.Range("K131").Value = 185952.528
$ws.Range("L131").Value = 13856.4198
$ws.Range("M131").Value = -180912.528
$ws.Range("N131").Value = -23936.4198
$ws.Range("H132").Value = 1798.25
$ws.Range("J132").Value = 1799.6666
$ws.Range("L132").Value = 16196.9994
$ws.Range("N132").Value = -21256.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 2778123
$ws.Range("I2").Value = 6250140
$ws.Range("J2").Value = 509.5
$ws.Range("K2").Value = 6250140
$ws.Range("L2").Value = 509.5
$ws.Range("M2").Value = -6250027
$ws.Range("N2").Value = -735.5
$ws.Range("H97").Value = 2370.818
$ws.Range("I97").Value = 2287.1538
$ws.Range("J97").Value = 2491.6667
$ws.Range("K97").Value = 2287.1538
$ws.Range("L97").Value = 2491.6667
$ws.Range("M97").Value = -1791.1538
$ws.Range("N97").Value = -3483.6667
$ws.Range("H126").Value = 6798.385
$ws.Range("I126").Value = 7230
$ws.Range("K126").Value = 21690
$ws.Range("M126").Value = -19220
$ws.Range("H132").Value = 3463.1482
$ws.Range("I132").Value = 3595.8125
$ws.Range("J132").Value = 3270.182
$ws.Range("K132").Value = 10787.4375
$ws.Range("L132").Value = 9810.545999999998
$ws.Range("M132").Value = -8257.4375
$ws.Range("N132").Value = -14870.546
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").ClearContents()
$ws.Range("N138").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 883.4706
$ws.Range("J46").Value = 2246.7144
$ws.Range("L46").Value = 2246.7144
$ws.Range("N46").Value = -2622.7144
$ws.Range("H68").Value = 3379.6667
$ws.Range("I68").Value = 2144.5
$ws.Range("J68").Value = 5850
$ws.Range("K68").Value = 2144.5
$ws.Range("L68").Value = 5850
$ws.Range("M68").Value = -1395.5
$ws.Range("N68").Value = -7348
$ws.Range("H71").Value = 3379.6667
$ws.Range("I71").Value = 2144.5
$ws.Range("J71").Value = 5850
$ws.Range("K71").Value = 10722.5
$ws.Range("L71").Value = 29250
$ws.Range("M71").Value = -6978.5
$ws.Range("N71").Value = -36738
$ws.Range("H81").Value = 79500
$ws.Range("J81").Value = 70000
$ws.Range("L81").Value = 70000
$ws.Range("N81").Value = -71996
$ws.Range("H84").Value = 79500
$ws.Range("J84").Value = 70000
$ws.Range("L84").Value = 210000
$ws.Range("N84").Value = -219984
$ws.Range("H93").Value = 1697.091
$ws.Range("J93").Value = 2123.842
$ws.Range("L93").Value = 2123.842
$ws.Range("N93").Value = -4619.842000000001
$ws.Range("H136").Value = 2490.14
$ws.Range("I136").Value = 1908.7567
$ws.Range("J136").Value = 4144.846
$ws.Range("K136").Value = 5726.2701
$ws.Range("L136").Value = 12434.538
$ws.Range("M136").Value = -3176.2701
$ws.Range("N136").Value = -17534.538

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 43999.668
$ws.Range("I70").Value = 30000
$ws.Range("K70").Value = 30000
$ws.Range("M70").Value = -29685
$ws.Range("H73").Value = 43999.668
$ws.Range("I73").Value = 30000
$ws.Range("K73").Value = 30000
$ws.Range("M73").Value = -28908
$ws.Range("H76").Value = 54500
$ws.Range("J76").Value = 51000
$ws.Range("L76").Value = 51000
$ws.Range("N76").Value = -51630
$ws.Range("H79").Value = 54500
$ws.Range("J79").Value = 51000
$ws.Range("L79").Value = 51000
$ws.Range("N79").Value = -53184
$ws.Range("H82").Value = 40742.43
$ws.Range("J82").Value = 40742.43
$ws.Range("L82").Value = 40742.43
$ws.Range("N82").Value = -41508.43
$ws.Range("H85").Value = 40742.43
$ws.Range("J85").Value = 40742.43
$ws.Range("L85").Value = 40742.43
$ws.Range("N85").Value = -43394.43
